# Apply the "assays version 1" edit:
#  - Insert two new leading columns ("version", "description") into the
#    "Export as TSV" sheet, shifting all existing columns two places right.
#  - Re-create the header-row cell comments at their shifted positions and
#    add comments for the two new columns.
#  - Add a data-validation rule for the new "version" column, sourced from
#    a new "version list" sheet.
#  - Insert the new "version list" sheet (containing the single value "1")
#    right after "Export as TSV".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export as TSV")

function Get-ColLetter {
    param([int]$colNum)
    $letter = ""
    while ($colNum -gt 0) {
        $rem = ($colNum - 1) % 26
        $letter = [char](65 + $rem) + $letter
        $colNum = [int](($colNum - $rem - 1) / 26)
    }
    return $letter
}

# ---------------------------------------------------------------------
# 1. Capture the existing header-row comments (column letter -> text)
#    before we shuffle anything around.
# ---------------------------------------------------------------------
$oldComments = @{}
foreach ($cmt in $ws.Comments) {
    $addr = $cmt.Parent.Address($false, $false)
    $oldComments[$addr] = $cmt.Text()
}

# Remove all existing comments; they'll be re-added in their new spots.
while ($ws.Comments.Count -gt 0) {
    $ws.Comments.Item(1).Delete()
}

# ---------------------------------------------------------------------
# 2. Insert two new columns at the front (A:B). Everything that used to
#    live in column A now lives in column C, etc. Data validations tied
#    to sqref ranges move automatically with the insert.
# ---------------------------------------------------------------------
$ws.Columns("A:B").Insert()

# ---------------------------------------------------------------------
# 3. Populate + style the two new header cells to match the existing
#    bold / centered / wrapped header look.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "version"
$ws.Cells.Item(1, 2).Value = "description"

$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Re-create the header comments, shifted two columns to the right,
#    plus two new ones for "version" and "description".
# ---------------------------------------------------------------------
$ws.Range("A1").AddComment("Version of the schema to use when validating this metadata.")
$ws.Range("B1").AddComment("Free-text description of this assay.")

foreach ($key in $oldComments.Keys) {
    # $key looks like "A1", "B1", ... "X1" -- pull the column index out,
    # shift it by two, and re-add the comment at the new address.
    $colIndex = $ws.Range($key).Column
    $newColIndex = $colIndex + 2
    $newColLetters = Get-ColLetter $newColIndex
    $targetAddr = "${newColLetters}1"
    $ws.Range($targetAddr).AddComment($oldComments[$key])
}

# ---------------------------------------------------------------------
# 5. Add the data-validation rule for the new "version" column.
# ---------------------------------------------------------------------
$versionRange = $ws.Range("A2:A1048576")
$versionRange.Validation.Add(
    [Microsoft.Office.Interop.Excel.XlDVType]::xlValidateList,
    [Microsoft.Office.Interop.Excel.XlFormatConditionOperator]::xlBetween,
    [System.Reflection.Missing]::Value,
    "='version list'!`$A`$1:`$A`$1"
)
$versionRange.Validation.ErrorTitle = "Value must come from list"
$versionRange.Validation.ErrorMessage = "Value must be one of: 1."
$versionRange.Validation.ShowInput = $true
$versionRange.Validation.ShowError = $true
$versionRange.Validation.IgnoreBlank = $true

# ---------------------------------------------------------------------
# 6. Insert the new "version list" sheet right after "Export as TSV".
# ---------------------------------------------------------------------
$versionSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$versionSheet.Name = "version list"
$versionSheet.Range("A1").Value = "1"
